# all_features.xlsx — apply cell-level formatting to the demo rows that
# showcase each supported style feature (border, color, fill, fontcolor,
# fontfamily, fontsize, halign, numberformat, text traits, valign).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlEdge* / xlVAlign* / xlHAlign* / xlLineStyle constants (no RGB() helper
# here, so colors are pre-computed decimal BGR-ish longs = R + G*256 + B*65536).
$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10

$xlContinuous = 1
$xlDash       = -4115

$xlHAlignLeft   = -4131
$xlHAlignCenter = -4108
$xlHAlignRight  = -4152

$xlVAlignTop    = -4160
$xlVAlignCenter = -4108
$xlVAlignBottom = -4107

$blueColor  = 16711680   # FF0000FF -> R=0x00 G=0x00 B=0xFF
$abcColor   = 16764091   # (A)BBCCFF -> R=0xBB G=0xCC B=0xFF (closest reachable w/o alpha)

# ---- border ----------------------------------------------------------
$ws.Range("B2").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous          # border=top
$ws.Range("B3").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous        # border=right
$ws.Range("B4").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous       # border=bottom
$ws.Range("B5").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous         # border=left
$ws.Range("B6").Borders.LineStyle = $xlContinuous                          # border=all

# border + borderstyle
$ws.Range("B9").Borders.LineStyle = $xlDash                                # border=all/borderstyle=dashed
$ws.Range("B10").Borders.LineStyle = $xlContinuous                         # border=all/borderstyle=solid

# ---- color (cell background fill) ------------------------------------
$ws.Range("B13").Interior.Color = $blueColor                               # color=FF0000
$ws.Range("B14").Interior.Color = $abcColor                                # color=ABC

# ---- fontcolor ---------------------------------------------------------
$ws.Range("B22").Font.Color = $blueColor                                   # fontcolor=FF0000
$ws.Range("B23").Font.Color = $abcColor                                    # fontcolor=ABC

# ---- fontfamily ----------------------------------------------------------
$ws.Range("B26").Font.Name = "Comic Sans Ms"                               # fontfamily='Comic Sans MS'
$ws.Range("B27").Font.Name = "Helvetica"                                   # fontfamily='Helvetica'

# ---- fontsize ----------------------------------------------------------
$ws.Range("B30").Font.Size = 20                                            # fontsize=20
$ws.Range("B31").Font.Size = 4                                             # fontsize=4

# ---- halign --------------------------------------------------------------
$ws.Range("B34").HorizontalAlignment = $xlHAlignLeft                       # halign=left
$ws.Range("B35").HorizontalAlignment = $xlHAlignCenter                     # halign=center
$ws.Range("B36").HorizontalAlignment = $xlHAlignRight                      # halign=right

# ---- numberformat --------------------------------------------------------
$ws.Range("C42").NumberFormat = "$#,##0_-"                                 # numberformat=currency
$ws.Range("C43").NumberFormat = "yyyy-mm-dd"                               # numberformat=date
$ws.Range("C44").NumberFormat = "d/m/yy h:mm"                              # numberformat=datetime
$ws.Range("C45").NumberFormat = "0"                                        # numberformat=number
$ws.Range("C46").NumberFormat = "0%"                                       # numberformat=percent
$ws.Range("C47").NumberFormat = "@"                                        # numberformat=text
$ws.Range("C48").NumberFormat = "h:mm AM/PM"                               # numberformat=time
$ws.Range("C49").NumberFormat = "0.00E+00"                                 # numberformat=scientific

# ---- text traits -----------------------------------------------------
$ws.Range("B52").Font.Bold = $true                                         # text=bold
$ws.Range("B53").Font.Italic = $true                                       # text=italic
$ws.Range("B54").Font.Underline = $true                                    # text=underline
$ws.Range("B55").Font.Strikethrough = $true                                # text=strikethrough

# ---- valign --------------------------------------------------------------
$ws.Range("B58").VerticalAlignment = $xlVAlignTop                          # valign=top
$ws.Range("B59").VerticalAlignment = $xlVAlignCenter                       # valign=center
$ws.Range("B60").VerticalAlignment = $xlVAlignBottom                       # valign=bottom
